{"js": "// Remove the \"Property\" / \"Value\" header row from the specification\n// table in the OVERVIEW section (the table whose data rows reference\n// the Jinja2 variable `overview_specifications_table[N].value`).\n//\n// The document contains several two-column tables that share the same\n// visual style (bold header cells), so we locate the right one by\n// inspecting each table's first row text instead of assuming a fixed\n// table index.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows.items;\n  if (rows.length === 0) {\n    continue;\n  }\n\n  const headerRow = rows[0];\n  headerRow.cells.load(\"items\");\n  await context.sync();\n\n  const cells = headerRow.cells.items;\n  for (const cell of cells) {\n    cell.body.load(\"text\");\n  }\n  await context.sync();\n\n  const cellTexts = cells.map((c) => c.body.text.trim());\n\n  // Only remove the header row from the table that actually has data\n  // rows below it (i.e. the specifications table, not the empty\n  // header-only tables used elsewhere in the template).\n  if (\n    cellTexts.length === 2 &&\n    cellTexts[0] === \"Property\" &&\n    cellTexts[1] === \"Value\" &&\n    rows.length > 1\n  ) {\n    headerRow.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Property\" / \"Value\" header row from the specification\n# table in the OVERVIEW section (the table whose data rows reference\n# the Jinja2 variable `overview_specifications_table[N].value`).\n#\n# The document contains several two-column tables that share the same\n# visual style (bold header cells), so we locate the right one by\n# inspecting each table's first row text instead of assuming a fixed\n# table index.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Tables.Count; $i -ge 1; $i--) {\n    $table = $d.Tables.Item($i)\n\n    if ($table.Rows.Count -le 1) {\n        continue\n    }\n\n    $headerRow = $table.Rows.Item(1)\n    if ($headerRow.Cells.Count -ne 2) {\n        continue\n    }\n\n    # Cell ranges include trailing cell-mark characters (CR + BEL), so\n    # strip those (not just whitespace) before comparing text.\n    $cell1 = $headerRow.Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7).Trim()\n    $cell2 = $headerRow.Cells.Item(2).Range.Text.TrimEnd([char]13, [char]7).Trim()\n\n    # Only remove the header row from the table that actually has data\n    # rows below it (i.e. the specifications table, not the empty\n    # header-only tables used elsewhere in the template).\n    if ($cell1 -eq \"Property\" -and $cell2 -eq \"Value\") {\n        $headerRow.Delete()\n    }\n}\n"}
